$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row of data appended after the last existing row (row 38 -> row 39)
$row = 39

# Column A holds a date string like "MM/DD/YYYY" stored as literal text
# (matching the pattern already used for recent rows, e.g. A38 = "07/04/2025").
# Excel's COM layer auto-detects strings that look like dates and silently
# converts them into date serial numbers, so we briefly force the cell to a
# text format while assigning the value, then restore a plain "Normal" style
# so no stray number formatting is left behind on the cell.
$cellA = $ws.Cells.Item($row, 1)
$cellA.NumberFormat = "@"
$cellA.Value = "07/06/2025"
$cellA.Style = "Normal"

$ws.Cells.Item($row, 2).Value = 131.3220000000001
$ws.Cells.Item($row, 3).Value = 0.07614870318758465
$ws.Cells.Item($row, 4).Value = 10
